$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting old D:K to F:M
$ws.Range("D:E").Insert()

# Copy number formatting from column F (old column D) into new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new Q4 2018 (D) and Q3 2018 (E) columns with reported figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 163400
$ws.Range("E8").Value = 92800
$ws.Range("D9").Value = 122900
$ws.Range("E9").Value = 68000
$ws.Range("D10").Value = 40500
$ws.Range("E10").Value = 24800
$ws.Range("D12").Value = 29700
$ws.Range("E12").Value = 26200
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 177100
$ws.Range("E17").Value = 116000
$ws.Range("D18").Value = -13700
$ws.Range("E18").Value = -23200
$ws.Range("D20").Value = 1500
$ws.Range("E20").Value = 1100
$ws.Range("D21").Value = -8300
$ws.Range("E21").Value = -18300
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -12300
$ws.Range("E23").Value = -22200
$ws.Range("D24").Value = 1000
$ws.Range("E24").Value = 200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -13300
$ws.Range("E26").Value = -22400
$ws.Range("D27").Value = -13300
$ws.Range("E27").Value = -22400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1500
$ws.Range("E32").Value = -1100
$ws.Range("D33").Value = -13300
$ws.Range("E33").Value = -22400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -13300
$ws.Range("E35").Value = -22400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 228400
$ws.Range("E41").Value = 166800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 87800
$ws.Range("E43").Value = 75000
$ws.Range("D44").Value = 80400
$ws.Range("E44").Value = 141600
$ws.Range("D45").Value = 23600
$ws.Range("E45").Value = 22500
$ws.Range("D46").Value = 420200
$ws.Range("E46").Value = 405900
$ws.Range("D47").Value = 9600
$ws.Range("E47").Value = 13000
$ws.Range("D48").Value = 35700
$ws.Range("E48").Value = 36300
$ws.Range("D49").Value = 17400
$ws.Range("E49").Value = 17700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 34200
$ws.Range("E52").Value = 33500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 517100
$ws.Range("E54").Value = 506400
$ws.Range("D57").Value = 32800
$ws.Range("E57").Value = 29000
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 96500
$ws.Range("E59").Value = 82800
$ws.Range("D60").Value = 129300
$ws.Range("E60").Value = 111700
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 44400
$ws.Range("E62").Value = 44400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 173800
$ws.Range("E66").Value = 156100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -306900
$ws.Range("E72").Value = -293500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 343400
$ws.Range("E76").Value = 350300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -13300
$ws.Range("E81").Value = -22400
$ws.Range("D83").Value = 3900
$ws.Range("E83").Value = 3900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 61400
$ws.Range("E89").Value = 40600
$ws.Range("D91").Value = -1500
$ws.Range("E91").Value = -500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -1500
$ws.Range("E94").Value = -100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 500
$ws.Range("E100").Value = -400
$ws.Range("D101").Value = 1200
$ws.Range("E101").Value = -400
$ws.Range("D102").Value = 61600
$ws.Range("E102").Value = 39700
